{"js": "// Append a new \"SourceCode\" paragraph at the end of the document body that\n// pulls the month out of the date column, mirroring the other R source-code\n// chunks already in this R-Markdown-derived report.\n//\n// Baltimore$month <- format (as.Date(Baltimore$date, format = \"%y/%m/%d\"), \"%m\")\n//\n// Each syntax-highlighted token is its own run with the matching Pandoc\n// \"Tok\" character style, exactly as the existing code chunks in the document\n// already do.\nconst runs = [\n  { text: \"Baltimore\", style: \"NormalTok\" },\n  { text: \"$\", style: \"SpecialCharTok\" },\n  { text: \"month \", style: \"NormalTok\" },\n  { text: \"<-\", style: \"OtherTok\" },\n  { text: \" \", style: \"NormalTok\" },\n  { text: \"format\", style: \"FunctionTok\" },\n  { text: \" (\", style: \"NormalTok\" },\n  { text: \"as.Date\", style: \"FunctionTok\" },\n  { text: \"(Baltimore\", style: \"NormalTok\" },\n  { text: \"$\", style: \"SpecialCharTok\" },\n  { text: \"date, \", style: \"NormalTok\" },\n  { text: \"format =\", style: \"AttributeTok\" },\n  { text: \" \", style: \"NormalTok\" },\n  { text: \"\\\"%y/%m/%d\\\"\", style: \"StringTok\" },\n  { text: \"), \", style: \"NormalTok\" },\n  { text: \"\\\"%m\\\"\", style: \"StringTok\" },\n  { text: \")\", style: \"NormalTok\" },\n];\n\nconst body = context.document.body;\n\n// New paragraph at the very end of the body, styled like the other code\n// blocks (\"Source Code\" is the display name of the SourceCode style).\nconst paragraph = body.insertParagraph(\"\", \"End\");\nparagraph.styleBuiltIn = Word.BuiltInStyleName.other;\nparagraph.style = \"Source Code\";\nawait context.sync();\n\n// Insert the first run's text directly into the (currently empty) new\n// paragraph, then keep chaining further runs after it, restyling each\n// inserted span with its own character style.\nlet range = paragraph.insertText(runs[0].text, \"End\");\nrange.style = runs[0].style;\nawait context.sync();\n\nfor (let i = 1; i < runs.length; i++) {\n  range = range.insertText(runs[i].text, \"After\");\n  range.style = runs[i].style;\n  await context.sync();\n}\n", "ps1": "# Append a new \"SourceCode\" paragraph at the end of the document body that\n# pulls the month out of the date column, mirroring the other R source-code\n# chunks already in this R-Markdown-derived report.\n#\n# Baltimore$month <- format (as.Date(Baltimore$date, format = \"%y/%m/%d\"), \"%m\")\n#\n# Each syntax-highlighted token is its own run with the matching Pandoc\n# \"Tok\" character style, exactly as the existing code chunks in the document\n# already do.\n$d = $word.ActiveDocument\n\n$runs = @(\n    @{ Text = 'Baltimore'; Style = 'NormalTok' }\n    @{ Text = '$'; Style = 'SpecialCharTok' }\n    @{ Text = 'month '; Style = 'NormalTok' }\n    @{ Text = '<-'; Style = 'OtherTok' }\n    @{ Text = ' '; Style = 'NormalTok' }\n    @{ Text = 'format'; Style = 'FunctionTok' }\n    @{ Text = ' ('; Style = 'NormalTok' }\n    @{ Text = 'as.Date'; Style = 'FunctionTok' }\n    @{ Text = '(Baltimore'; Style = 'NormalTok' }\n    @{ Text = '$'; Style = 'SpecialCharTok' }\n    @{ Text = 'date, '; Style = 'NormalTok' }\n    @{ Text = 'format ='; Style = 'AttributeTok' }\n    @{ Text = ' '; Style = 'NormalTok' }\n    @{ Text = '\"%y/%m/%d\"'; Style = 'StringTok' }\n    @{ Text = '), '; Style = 'NormalTok' }\n    @{ Text = '\"%m\"'; Style = 'StringTok' }\n    @{ Text = ')'; Style = 'NormalTok' }\n)\n\n# Add a brand-new paragraph after everything else in the body.\n$endRange = $d.Content\n$endRange.Collapse(0)\n[void]$endRange.InsertParagraphAfter()\n\n# Grab that freshly-minted (still empty) last paragraph and give it the\n# \"Source Code\" paragraph style used by every other code chunk.\n$p = $d.Paragraphs.Last\n$p.Style = \"Source Code\"\n$pr = $p.Range\n\n# Seed the paragraph with the first token, then trim the trailing paragraph\n# mark off the range before stamping the run-level character style so the\n# style lands on the run, not the paragraph.\n$pr.Text = $runs[0].Text\n[void]$pr.MoveEnd(1, -1)\n$pr.Style = $runs[0].Style\n\n# Every following token gets inserted right after the previous run's end and\n# re-styled the same way.\n$insertAt = $pr.End\nfor ($i = 1; $i -lt $runs.Count; $i++) {\n    $run = $d.Range($insertAt, $insertAt)\n    [void]$run.InsertAfter($runs[$i].Text)\n    $run.Style = $runs[$i].Style\n    $insertAt = $run.End\n}\n"}
